$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = "B0C9THQ36P"
$ws.Range("B11").Value = "XVX Retro 75% Gaming Keyboard with OLED Display&Knob, M87 Pro Bluetooth 5.1/2.4GHz /USB-C Wireless Mechanical Keyboard with Hot-Swappable Custom Switch, Compact TKL Gamer RGB Keyboard, PBT Keycaps"
$ws.Range("C11").Value = 392.4
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "02/03/2024"
